# Update automática del mapa (2025-08-11 11:27:20)
# The oldest pending case (row 55: Caso 6104 - PINTO 4677) has been
# resolved/removed, so the whole row is deleted and every row below it
# shifts up by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(55).Delete()
